$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 69, shifting existing rows 69:88 down to 70:89
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly data point
$ws.Range("A69").Value = 7
$ws.Range("B69").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C69").Value = "Ñuble"
$ws.Range("D69").Value = 45211
$ws.Range("D69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E69").Value = 16
$ws.Range("F69").Value = 100112026
$ws.Range("G69").Value = "Haba"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 60
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = 12000
$ws.Range("N69").Value = "$/saco 25 kilos"
$ws.Range("O69").Value = "Provincia de Diguillín"
$ws.Range("P69").Value = 480
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
